# Overall_Rebate_Efficiency.xlsx — "Add files via upload"
#
# The PSA_LOLO sheet is reshaped from a single header row of two named
# columns (psa_lolo_20 / psa_lolo_40) holding one data row each, into a
# two-column (label, value) table with a single "PSALOLO" header and two
# data rows (20 -> 30197, 40 -> 15279).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PSA_LOLO")

# Wipe the old 2x2 block (headers + single data row) completely, content
# and formatting both, so stale styles don't linger on cells that should
# end up blank.
$ws.Range("A1:C3").Clear()

# New header: just one label now, sitting in B1 (A1 stays blank), with a
# blank-but-styled neighbor in C1 matching the sheet's small header font.
$ws.Range("B1").Value = "PSALOLO"
$ws.Range("B1").Font.Size = 9
$ws.Range("C1").Font.Size = 9

# New data rows: the old "_20"/"_40" column-name suffixes become the
# row labels in column A, paired with the two original data values.
$ws.Range("A2").Value = 20
$ws.Range("B2").Value = 30197
$ws.Range("A3").Value = 40
$ws.Range("B3").Value = 15279

# Reset the cursor back to the top of the sheet now that the rebuild is
# done (rather than leaving it parked on the stale pre-edit selection).
$ws.Range("A1").Select() | Out-Null
